# fix: fix a conflict
#
# 1. Widen the saved window (bookViews/workbookView windowWidth 18288 -> 23040).
# 2. Turn on iterative calculation (calcPr iterate="1" iterateCount="100" iterateDelta="0.001").
# 3. Rename the header in A1 from "num" to "id".
# 4. Move the worksheet's remembered selection from E5 to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Resize the workbook window ---------------------------------------
$excel.ActiveWindow.Width = 23040

# --- 2. Enable iterative calculation (100 iterations, max change 0.001) ---
$excel.Iteration = $true
$excel.MaxIterations = 100
$excel.MaxChange = 0.001

# --- 3. Fix the conflicting header text in A1 ------------------------------
$ws.Range("A1").Value = "id"

# --- 4. Update the remembered selection to A2 ------------------------------
[void]$ws.Range("A2").Select()
